$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add the new "ODI Bowling Extra" worksheet as the last tab.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Bowling Extra"

# Headers
$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "MAIDEN_OVERS"
$extra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"
$extra.Range("A1:C1").Font.Bold = $true
$extra.Range("A1:C1").HorizontalAlignment = -4108
$extra.Range("A1:C1").VerticalAlignment = -4160
$extra.Range("A1:C1").Borders.LineStyle = 1

# Data rows (MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL)
$rows = @(
    @("3650", "0", ""),
    @("3916", "", ""),
    @("3917", "0", "20.00%"),
    @("3918", "0", ""),
    @("3938", "", ""),
    @("3941", "", ""),
    @("4009", "0", ""),
    @("4129", "0", "10.00%"),
    @("4131", "0", ""),
    @("4140", "0", ""),
    @("4198", "1", "10.00%"),
    @("4203", "", ""),
    @("4257", "0", ""),
    @("4290", "0", ""),
    @("4301", "0", "10.00%"),
    @("4315", "0", ""),
    @("4326", "", ""),
    @("4332", "", ""),
    @("4335", "0", ""),
    @("4538", "", "")
)

$r = 2
foreach ($row in $rows) {
    $extra.Cells.Item($r, 1).NumberFormat = "@"
    $extra.Cells.Item($r, 1).Value = $row[0]
    $extra.Cells.Item($r, 2).NumberFormat = "@"
    $extra.Cells.Item($r, 2).Value = $row[1]
    $extra.Cells.Item($r, 3).NumberFormat = "@"
    $extra.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# ---------------------------------------------------------------------------
# 2) Clean up the "ODI Batting Extra" sheet: drop the empty placeholder
#    cells in B10:E10 and B14:E14 (rows for match codes 4525 and 4538,
#    which have no batting-position/num4/num6/percent data).
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$battingExtra.Range("B10:E10").ClearContents()
$battingExtra.Range("B14:E14").ClearContents()
